# Adds "Ecology"/"Perch Data" ecomorph-code columns (E:F) to the species
# list, fixes two mis-tagged Region values, and turns the A1:F201 range into
# a filtered table (matches commit "add code for l1ou").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MainlandAnole_SpeciesList")

# --- header row -----------------------------------------------------------
$ws.Range("E1").Value = "Ecology"
$ws.Range("F1").Value = "Perch Data"

# --- per-species ecology code (E) + perch-data flag (F) -------------------
# Each tuple is (row, EcologyCode, PerchDataFlag). Row 1 (header) handled above.
$rows = @(
    @(2, 'TG', 'Y'),
    @(3, 'Tw', 'Y'),
    @(5, 'TC', 'Y'),
    @(6, 'TC', 'Y'),
    @(7, 'TG', 'Y'),
    @(10, 'GB', 'Y'),
    @(12, 'Tw', 'Y'),
    @(15, 'U', 'Y'),
    @(18, 'GB', 'Y'),
    @(19, 'CG', 'Y'),
    @(20, 'U', 'Y'),
    @(22, 'U', 'Y'),
    @(31, 'T', 'Y'),
    @(34, 'T', 'Y'),
    @(37, 'TC', 'Y'),
    @(38, 'U', 'Y'),
    @(41, 'TC', 'Y'),
    @(43, 'TG', 'Y'),
    @(44, 'TG', 'Y'),
    @(45, 'TG', 'Y'),
    @(49, 'CG', 'Y'),
    @(52, 'T', 'Y'),
    @(56, 'CG', 'Y'),
    @(60, 'TC', 'Y'),
    @(64, 'CG', 'Y'),
    @(65, 'Tw', 'Y'),
    @(67, 'TC', 'Y'),
    @(69, 'TG', 'Y'),
    @(70, 'GB', 'Y'),
    @(74, 'TG ', 'Y'),
    @(76, 'Twig', 'Y'),
    @(77, 'U', 'Y'),
    @(79, 'TG', 'Y'),
    @(81, 'GB', 'Y'),
    @(83, 'GB', 'Y'),
    @(86, 'Y', 'GB'),
    @(91, 'U', 'Y'),
    @(93, 'T', 'Y'),
    @(94, 'U', 'Y'),
    @(95, 'CG', 'Y'),
    @(101, 'TG', 'Y'),
    @(105, 'TC', 'Y'),
    @(110, 'TG', 'Y'),
    @(113, 'TG', 'Y'),
    @(118, 'CG', 'Y'),
    @(121, 'GB', 'Y'),
    @(124, 'TC', 'Y'),
    @(125, 'GB', 'Y'),
    @(138, 'U', 'Y'),
    @(139, 'TC', 'Y'),
    @(140, 'U', 'Y'),
    @(141, 'GB', 'Y'),
    @(142, 'U', 'Y'),
    @(146, 'TG', 'Y'),
    @(147, 'GB', 'Y'),
    @(148, 'CG', 'Y'),
    @(156, 'TG', 'Y'),
    @(159, 'U', 'Y'),
    @(160, 'TG', 'Y'),
    @(162, 'GB', 'Y'),
    @(164, 'Twig', 'Y'),
    @(168, 'TC', 'Y'),
    @(178, 'Twig', 'Y'),
    @(179, 'GB', 'Y'),
    @(183, 'U', 'Y'),
    @(191, 'U', 'Y'),
    @(192, 'CG', 'Y'),
    @(200, 'GB', 'Y')
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 5).Value = $row[1]
    $ws.Cells.Item($r, 6).Value = $row[2]
}

# --- fix Region for guamuhaya (190) and vermiculatus (191): Mainland -> Caribbean
$ws.Range("D190").Value = "Caribbean"
$ws.Range("D191").Value = "Caribbean"

# --- turn the used range into a filtered table -----------------------------
[void]$ws.Range("A1:F201").AutoFilter()

$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=MainlandAnole_SpeciesList!`$A`$1:`$F`$201")
$filterName.Visible = $false

# --- move the selection (matches the saved cursor position in the diff) ----
[void]$ws.Range("A130").Select()
